$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two id values in column A (rows 2 and 3) to new random strings
$ws.Range("A2").Value = "dmxQwKymKD3FrUgJHgCr"
$ws.Range("A3").Value = "wYfhX0ordSBl1agNeVgm"
